$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the
#    first paragraph (the "Play Buffalo Mania Slot Game for Free"
#    Heading1). We build it by literally moving the duplicate bold
#    paragraph that currently sits near the end of the document
#    (this preserves its original run layout, including the leading
#    empty run) and then edit its text in place.
# ------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$pBoldDup = $d.Paragraphs.Item($n - 1)
$boldFormattedText = $pBoldDup.Range.FormattedText

$pMeta = $d.Paragraphs.Item(2)
$pMeta.Style = "Normal"
$pMeta.Range.FormattedText = $boldFormattedText

# Turn the copied "Play Buffalo Mania Slot Game for Free" bold run
# into "Meta description".
$pMetaRange = $pMeta.Range
$oldBoldText = "Play Buffalo Mania Slot Game for Free"
$boldRun = $d.Range($pMetaRange.Start, $pMetaRange.Start + $oldBoldText.Length)
$boldRun.Text = "Meta description"

# Append the (non-bold) rest of the meta description text.
$pMetaRange = $pMeta.Range
$insertionPoint = $d.Range($pMetaRange.End - 1, $pMetaRange.End - 1)
$metaRest = ": Read our review of Buffalo Mania, an online slot game featuring 720 possible ways to win and three progressive jackpots. Play now for free!"
$insertionPoint.InsertAfter($metaRest)

# ------------------------------------------------------------------
# 2) Remove the now-redundant bold "Play Buffalo Mania Slot Game for
#    Free" paragraph that used to live near the end of the document.
# ------------------------------------------------------------------
$n2 = $d.Paragraphs.Count
$pBoldDup2 = $d.Paragraphs.Item($n2 - 1)
$pBoldDup2.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the italic "Read our review..." text (now the very last
#    paragraph) with the new image-prompt copy, keeping the italic
#    formatting and avoiding any smart-quote autocorrection.
# ------------------------------------------------------------------
$n3 = $d.Paragraphs.Count
$pImg = $d.Paragraphs.Item($n3)
$pImgRange = $pImg.Range
$oldText = $d.Range($pImgRange.Start, $pImgRange.End - 1)
$oldText.Delete()

$newText = 'Create a feature image fitting the game "Buffalo Mania" with a happy Maya warrior wearing glasses in cartoon style. The image should showcase the warrior riding on a majestic bison through the North American prairies with the game''s logo in bold letters. The background should highlight the Yellowstone''s snow-capped peaks and the wind with flutes. Use bright and vibrant colors to appeal to the audience.'

$pImg2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$insPoint = $d.Range($pImg2.Range.Start, $pImg2.Range.Start)
$insPoint.InsertAfter($newText)

$pImg3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pImg3Range = $pImg3.Range
$finalTextRange = $d.Range($pImg3Range.Start, $pImg3Range.End - 1)
$finalTextRange.Italic = 1
